{"js": "// Replace the 25 multiplication problems in the practice table with the\n// new values described by the commit diff. Each lookup value is unique in\n// the document, so a plain text search + replace is unambiguous.\nconst replacements = [\n  [\"964\u00d76=\", \"756\u00d79=\"],\n  [\"872\u00d76=\", \"928\u00d79=\"],\n  [\"880\u00d79=\", \"358\u00d75=\"],\n  [\"143\u00d77=\", \"681\u00d74=\"],\n  [\"639\u00d72=\", \"431\u00d79=\"],\n  [\"318\u00d78=\", \"253\u00d73=\"],\n  [\"792\u00d77=\", \"904\u00d77=\"],\n  [\"969\u00d79=\", \"801\u00d77=\"],\n  [\"633\u00d79=\", \"447\u00d75=\"],\n  [\"580\u00d79=\", \"920\u00d72=\"],\n  [\"683\u00d77=\", \"322\u00d74=\"],\n  [\"521\u00d72=\", \"118\u00d72=\"],\n  [\"869\u00d76=\", \"910\u00d79=\"],\n  [\"260\u00d74=\", \"948\u00d73=\"],\n  [\"616\u00d72=\", \"922\u00d74=\"],\n  [\"793\u00d79=\", \"740\u00d72=\"],\n  [\"950\u00d77=\", \"255\u00d73=\"],\n  [\"629\u00d76=\", \"874\u00d79=\"],\n  [\"501\u00d75=\", \"298\u00d74=\"],\n  [\"986\u00d77=\", \"169\u00d76=\"],\n  [\"595\u00d75=\", \"642\u00d77=\"],\n  [\"734\u00d74=\", \"356\u00d73=\"],\n  [\"797\u00d76=\", \"190\u00d76=\"],\n  [\"791\u00d76=\", \"485\u00d79=\"],\n  [\"277\u00d78=\", \"450\u00d77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 multiplication problems in the practice table with the\n# new values described by the commit diff. Each lookup value is unique in\n# the document, so Find/Replace is unambiguous.\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$pairs = @(\n  @{old=\"964\u00d76=\"; new=\"756\u00d79=\"},\n  @{old=\"872\u00d76=\"; new=\"928\u00d79=\"},\n  @{old=\"880\u00d79=\"; new=\"358\u00d75=\"},\n  @{old=\"143\u00d77=\"; new=\"681\u00d74=\"},\n  @{old=\"639\u00d72=\"; new=\"431\u00d79=\"},\n  @{old=\"318\u00d78=\"; new=\"253\u00d73=\"},\n  @{old=\"792\u00d77=\"; new=\"904\u00d77=\"},\n  @{old=\"969\u00d79=\"; new=\"801\u00d77=\"},\n  @{old=\"633\u00d79=\"; new=\"447\u00d75=\"},\n  @{old=\"580\u00d79=\"; new=\"920\u00d72=\"},\n  @{old=\"683\u00d77=\"; new=\"322\u00d74=\"},\n  @{old=\"521\u00d72=\"; new=\"118\u00d72=\"},\n  @{old=\"869\u00d76=\"; new=\"910\u00d79=\"},\n  @{old=\"260\u00d74=\"; new=\"948\u00d73=\"},\n  @{old=\"616\u00d72=\"; new=\"922\u00d74=\"},\n  @{old=\"793\u00d79=\"; new=\"740\u00d72=\"},\n  @{old=\"950\u00d77=\"; new=\"255\u00d73=\"},\n  @{old=\"629\u00d76=\"; new=\"874\u00d79=\"},\n  @{old=\"501\u00d75=\"; new=\"298\u00d74=\"},\n  @{old=\"986\u00d77=\"; new=\"169\u00d76=\"},\n  @{old=\"595\u00d75=\"; new=\"642\u00d77=\"},\n  @{old=\"734\u00d74=\"; new=\"356\u00d73=\"},\n  @{old=\"797\u00d76=\"; new=\"190\u00d76=\"},\n  @{old=\"791\u00d76=\"; new=\"485\u00d79=\"},\n  @{old=\"277\u00d78=\"; new=\"450\u00d77=\"}\n)\n\n$d = $word.ActiveDocument\n\nforeach ($p in $pairs) {\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($p.old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $p.new, $wdReplaceAll)\n}\n"}
